# Apply cryptos.xlsx update: refresh Price (D) and Volume(1h) (E) columns,
# and fix the swapped Polkadot/Chainlink rows (14 & 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.532.97'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.482.09'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = "'313.54"
$ws.Range('D5').ClearFormats()
$ws.Range('D6').Value = "'92.38"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.27%  '
$ws.Range('D7').Value = "'0.550"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('D10').Value = "'32.79"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('D11').Value = "'0.0793"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('D13').Value = '2.865.78'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'16.45"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +9.93%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = "'6.94"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '2.477.74'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').Value = "'0.778"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '41.576.86'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  +3.73%  '
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').Value = "'72.50"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +5.62%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = "'236.78"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = "'25.05"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = "'9.72"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').Value = "'36.02"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').Value = "'157.70"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.97%  '
$ws.Range('D32').Value = "'5.46"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').Value = "'2.57"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D34').Value = "'0.0758"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('D35').Value = "'17.59"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.77%  '
$ws.Range('E36').Value = '  -9.95%  '
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('E38').Value = '  -5.01%  '
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = '1.970.71'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D44').Value = "'19.05"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.50%  '
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('D47').Value = "'8.98"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('D48').Value = '2.722.84'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = "'98.09"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').Value = "'68.26"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('E51').Value = '  -2.26%  '
